$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "sno" -> "id", "number" -> "age". "name" stays the same.
$ws.Range("A1").Value = "id"
$ws.Range("C1").Value = "age"

# Move the active selection to A2 (previously C5)
$ws.Range("A2").Select()
